$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SynonymPredicateChange")

# Insert a new column before column D (shifts about_node..has_undo right by one)
$ws.Range("D1").EntireColumn.Insert()

# Set the header of the newly inserted column D
$ws.Range("D1").Value = "predicate"
